$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '71.015.95'
$ws.Range('E2').Value = '  +6.06%  '
$ws.Range('D3').Value = '3.664.27'
$ws.Range('E3').Value = '  +17.93%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '620.07'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +7.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '181.83'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.33%  '
$ws.Range('D7').Value = '3.663.82'
$ws.Range('E7').Value = '  +17.93%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('E9').Value = '  +5.62%  '
$ws.Range('E10').Value = '  +7.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.67'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +4.92%  '
$ws.Range('E12').Value = '  +7.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '40.34'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +11.48%  '
$ws.Range('E14').Value = '  +6.24%  '
$ws.Range('D15').Value = '4.272.42'
$ws.Range('E15').Value = '  +17.88%  '
$ws.Range('D16').Value = '70.999.73'
$ws.Range('E16').Value = '  +6.09%  '
$ws.Range('D17').Value = '3.667.18'
$ws.Range('E17').Value = '  +18.15%  '
$ws.Range('E18').Value = '  +2.06%  '
$ws.Range('E19').Value = '  +7.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '520.25'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +8.43%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.92'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.75%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.23'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +18.55%  '
$ws.Range('E23').Value = '  +7.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.53'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +13.16%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '88.56'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +5.87%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '13.50'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +7.31%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.08'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +9.74%  '
$ws.Range('E28').Value = '  -0.05%  '
$ws.Range('E29').Value = '  +11.17%  '
$ws.Range('E30').Value = '  +3.57%  '
$ws.Range('E31').Value = '  +11.99%  '
$ws.Range('E32').Value = '  +17.74%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '31.60'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +12.86%  '
$ws.Range('E34').Value = '  +4.57%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.12'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +9.39%  '
$ws.Range('E37').Value = '  +9.08%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.349'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +12.05%  '
$ws.Range('E39').Value = '  +9.55%  '
$ws.Range('E40').Value = '  +7.05%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '51.28'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +4.49%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '45.37'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -6.72%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '432.30'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +15.90%  '
$ws.Range('E44').Value = '  +6.01%  '
$ws.Range('D45').Value = '3.110.76'
$ws.Range('E45').Value = '  +11.09%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.81'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +4.01%  '
$ws.Range('E47').Value = '  +7.77%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '28.26'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +9.86%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '139.75'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +3.01%  '
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('E51').Value = '  +11.03%  '
